# JobMaterial.xlsx — rotate the Material rows.
#
# Net effect required by the change:
#   Row 3 <- (original) Row 5
#   Row 4 <- (original) Row 3
#   Row 5 <- (original) Row 4
#   Row 6 <- (original) Row 7
#   Row 7 <- (original) Row 6
#
# Columns B:H carry the data for each row (A is a repeated section label and
# is left untouched). We stage the five rows' original B:H content into a
# scratch area first (Copy captures a snapshot at call time, so staging
# avoids clobbering a row before we've read the value we still need from
# it), then paste the staged snapshots into their new homes using
# "paste values" so the destination cell keeps its own formatting/style and
# only the text content + cell type change - exactly like a normal
# Excel value-only paste.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

# --- Stage original row contents (rows 3-7) into scratch rows 100-104 ---
$ws.Range("B3:H3").Copy()
$ws.Range("B100:H100").PasteSpecial($xlPasteValues)

$ws.Range("B4:H4").Copy()
$ws.Range("B101:H101").PasteSpecial($xlPasteValues)

$ws.Range("B5:H5").Copy()
$ws.Range("B102:H102").PasteSpecial($xlPasteValues)

$ws.Range("B6:H6").Copy()
$ws.Range("B103:H103").PasteSpecial($xlPasteValues)

$ws.Range("B7:H7").Copy()
$ws.Range("B104:H104").PasteSpecial($xlPasteValues)

# --- Write the rotated content back from the scratch snapshots ---
# Row 3 <- original Row 5 (scratch 102)
$ws.Range("B102:H102").Copy()
$ws.Range("B3:H3").PasteSpecial($xlPasteValues)

# Row 4 <- original Row 3 (scratch 100)
$ws.Range("B100:H100").Copy()
$ws.Range("B4:H4").PasteSpecial($xlPasteValues)

# Row 5 <- original Row 4 (scratch 101)
$ws.Range("B101:H101").Copy()
$ws.Range("B5:H5").PasteSpecial($xlPasteValues)

# Row 6 <- original Row 7 (scratch 104)
$ws.Range("B104:H104").Copy()
$ws.Range("B6:H6").PasteSpecial($xlPasteValues)

# Row 7 <- original Row 6 (scratch 103)
$ws.Range("B103:H103").Copy()
$ws.Range("B7:H7").PasteSpecial($xlPasteValues)

# --- Clean up the scratch area ---
$ws.Range("B100:H104").Clear()
